# Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resultados (Acierto / Fallo) y profit asociados para cada fila actualizada.
# profit = cuota - 1 cuando el resultado es "Acierto"; profit = -1 cuando es "Fallo".
$updates = @(
    @{ Row = 139; Resultado = "Acierto"; Profit = 0.36 },
    @{ Row = 141; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 143; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 144; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 152; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 156; Resultado = "Acierto"; Profit = 2.5 },
    @{ Row = 157; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 160; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 162; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 163; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 164; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 166; Resultado = "Acierto"; Profit = 0.83 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Resultado
    $ws.Cells.Item($u.Row, 8).Value = $u.Profit
}

# event_id de la fila 167 estaba almacenado como texto; se corrige a numérico.
$ws.Cells.Item(167, 1).Value = 14689564
